$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63; this shifts existing rows 63:100 down to 64:101
$ws.Rows.Item(63).Insert()

# Populate the new row 63 with a new weekly record (same fixed fields as the
# surrounding rows, new date + price figures)
$ws.Cells.Item(63, 1).Value = 1
$ws.Cells.Item(63, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(63, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(63, 4).Value = 44813
$ws.Cells.Item(63, 5).Value = 15
$ws.Cells.Item(63, 6).Value = 100112038
$ws.Cells.Item(63, 7).Value = "Cebollín baby"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 250
$ws.Cells.Item(63, 11).Value = 1200
$ws.Cells.Item(63, 12).Value = 1500
$ws.Cells.Item(63, 13).Value = 1350
$ws.Cells.Item(63, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(63, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(63, 16).Value = 675
$ws.Cells.Item(63, 17).Value = 2
$ws.Cells.Item(63, 18).Value = "Hortaliza"
